$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed values in existing rows (rows 3-26) ---
$ws.Range("B3").Value = 0.7941176470588235
$ws.Range("C3").Value = 27
$ws.Range("D3").Value = 27
$ws.Range("H3").Value = 7
$ws.Range("J3").Value = "love"
$ws.Range("K3").Value = 0.9782608695652174
$ws.Range("L3").Value = 45
$ws.Range("M3").Value = 45
$ws.Range("Q3").Value = 1
$ws.Range("A4").Value = "crisis"
$ws.Range("B4").Value = 0.5993150684931506
$ws.Range("C4").Value = 175
$ws.Range("D4").Value = 175
$ws.Range("H4").Value = 117
$ws.Range("J4").Value = "interesting"
$ws.Range("K4").Value = 0.9393939393939394
$ws.Range("Q4").Value = 2
$ws.Range("A5").Value = "sc"
$ws.Range("B5").Value = 0.1746031746031746
$ws.Range("C5").Value = 33
$ws.Range("D5").Value = 33
$ws.Range("H5").Value = 156
$ws.Range("J5").Value = "best"
$ws.Range("K5").Value = 0.9322033898305084
$ws.Range("L5").Value = 55
$ws.Range("M5").Value = 55
$ws.Range("Q5").Value = 4
$ws.Range("B6").Value = 0.1744186046511628
$ws.Range("C6").Value = 90
$ws.Range("D6").Value = 90
$ws.Range("H6").Value = 426
$ws.Range("J6").Value = "great"
$ws.Range("K6").Value = 0.8660714285714286
$ws.Range("L6").Value = 97
$ws.Range("M6").Value = 97
$ws.Range("Q6").Value = 15
$ws.Range("J7").Value = "won"
$ws.Range("K7").Value = 0.8205128205128205
$ws.Range("L7").Value = 32
$ws.Range("M7").Value = 32
$ws.Range("Q7").Value = 7
$ws.Range("J8").Value = "thanks"
$ws.Range("K8").Value = 0.8170731707317073
$ws.Range("L8").Value = 67
$ws.Range("M8").Value = 67
$ws.Range("Q8").Value = 15
$ws.Range("K9").Value = 0.7916666666666666
$ws.Range("L9").Value = 95
$ws.Range("M9").Value = 95
$ws.Range("Q9").Value = 25
$ws.Range("K10").Value = 0.7890625
$ws.Range("L10").Value = 101
$ws.Range("M10").Value = 101
$ws.Range("Q10").Value = 27
$ws.Range("J11").Value = "special"
$ws.Range("K11").Value = 0.7777777777777778
$ws.Range("L11").Value = 28
$ws.Range("M11").Value = 28
$ws.Range("Q11").Value = 8
$ws.Range("J12").Value = "positive"
$ws.Range("K12").Value = 0.7758620689655172
$ws.Range("L12").Value = 45
$ws.Range("M12").Value = 45
$ws.Range("Q12").Value = 13
$ws.Range("J13").Value = "support"
$ws.Range("K13").Value = 0.7358490566037735
$ws.Range("L13").Value = 78
$ws.Range("M13").Value = 78
$ws.Range("Q13").Value = 28
$ws.Range("J14").Value = "confidence"
$ws.Range("K14").Value = 0.7222222222222222
$ws.Range("L14").Value = 26
$ws.Range("M14").Value = 26
$ws.Range("K15").Value = 0.7125
$ws.Range("L15").Value = 114
$ws.Range("M15").Value = 114
$ws.Range("Q15").Value = 46
$ws.Range("J16").Value = "safe"
$ws.Range("K16").Value = 0.7112676056338029
$ws.Range("L16").Value = 101
$ws.Range("M16").Value = 101
$ws.Range("Q16").Value = 41
$ws.Range("J17").Value = "safety"
$ws.Range("K17").Value = 0.6666666666666666
$ws.Range("L17").Value = 34
$ws.Range("M17").Value = 34
$ws.Range("Q17").Value = 17
$ws.Range("J18").Value = "better"
$ws.Range("K18").Value = 0.6349206349206349
$ws.Range("L18").Value = 40
$ws.Range("M18").Value = 40
$ws.Range("Q18").Value = 23
$ws.Range("J19").Value = "fresh"
$ws.Range("K19").Value = 0.5833333333333334
$ws.Range("L19").Value = 28
$ws.Range("M19").Value = 28
$ws.Range("Q19").Value = 20
$ws.Range("J20").Value = "well"
$ws.Range("K20").Value = 0.574468085106383
$ws.Range("L20").Value = 54
$ws.Range("M20").Value = 54
$ws.Range("Q20").Value = 40
$ws.Range("K21").Value = 0.56
$ws.Range("L21").Value = 28
$ws.Range("M21").Value = 28
$ws.Range("Q21").Value = 22
$ws.Range("K22").Value = 0.4986945169712794
$ws.Range("L22").Value = 191
$ws.Range("M22").Value = 191
$ws.Range("Q22").Value = 192
$ws.Range("K23").Value = 0.4705882352941176
$ws.Range("L23").Value = 160
$ws.Range("M23").Value = 160
$ws.Range("Q23").Value = 180
$ws.Range("J24").Value = "care"
$ws.Range("K24").Value = 0.449438202247191
$ws.Range("L24").Value = 40
$ws.Range("M24").Value = 40
$ws.Range("Q24").Value = 49
$ws.Range("J25").Value = "help"
$ws.Range("K25").Value = 0.4169491525423729
$ws.Range("L25").Value = 123
$ws.Range("M25").Value = 123
$ws.Range("Q25").Value = 172
$ws.Range("J26").Value = "protect"
$ws.Range("K26").Value = 0.410958904109589
$ws.Range("L26").Value = 30
$ws.Range("M26").Value = 30
$ws.Range("Q26").Value = 43

# --- Remove cells for shrunk A-table (rows 7-8, cols A-H) ---
$ws.Range("A7:H8").Clear()

# --- Add new row 27 (J27:Q27) ---
$ws.Range("J27").Value = "please"
$ws.Range("K27").Value = 0.3430962343096234
$ws.Range("L27").Value = 82
$ws.Range("M27").Value = 82
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 157

# --- Copy style from J26 (existing styled row) to J27 ---
$ws.Range("J26").Copy()
$ws.Range("J27").PasteSpecial(-4122)
